$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.824.72'
$ws.Range("E2").Value = '  -0.47%  '

$ws.Range("D3").Value = '1.886.07'
$ws.Range("E3").Value = '  -0.52%  '

$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '0.7497'
$r.Style = "Normal"
$ws.Range("E5").Value = '  -2.98%  '

$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '241.97'
$r.Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '

$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '1.0000'
$r.Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '

$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '0.3126'
$r.Style = "Normal"
$ws.Range("E8").Value = '  -0.50%  '

$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '25.25'
$r.Style = "Normal"
$ws.Range("E9").Value = '  -2.48%  '

$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.07109'
$r.Style = "Normal"
$ws.Range("E10").Value = '  -3.50%  '

$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.08505'
$r.Style = "Normal"
$ws.Range("E11").Value = '  +5.36%  '

$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.7588'
$r.Style = "Normal"
$ws.Range("E12").Value = '  -1.81%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.876.10'
$ws.Range("E13").Value = '  -1.40%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '5.362'
$r.Style = "Normal"
$ws.Range("E14").Value = '  -2.60%  '

$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '93.18'
$r.Style = "Normal"
$ws.Range("E15").Value = '  -1.28%  '

$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '6.124'
$r.Style = "Normal"
$ws.Range("E16").Value = '  -1.75%  '

$ws.Range("D17").Value = '29.605.10'
$ws.Range("E17").Value = '  -1.24%  '

$ws.Range("E18").Value = '  -2.33%  '

$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = '242.69'
$r.Style = "Normal"
$ws.Range("E19").Value = '  -2.05%  '

$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '0.000007833'
$r.Style = "Normal"
$ws.Range("E20").Value = '  -0.21%  '

$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '0.9994'
$r.Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '2.140.95'
$ws.Range("E22").Value = '  -1.19%  '

$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '7.941'
$r.Style = "Normal"
$ws.Range("E23").Value = '  -2.84%  '

$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '1.000'
$r.Style = "Normal"
$ws.Range("E24").Value = '  -0.08%  '

$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.1590'
$r.Style = "Normal"
$ws.Range("E25").Value = '  +0.80%  '

$ws.Range("E26").Value = '  -1.06%  '

$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = '162.95'
$r.Style = "Normal"
$ws.Range("E27").Value = '  -0.09%  '

$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = '18.69'
$r.Style = "Normal"
$ws.Range("E28").Value = '  -0.43%  '

$ws.Range("E29").Value = '  -0.37%  '

$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = '1.474'
$r.Style = "Normal"
$ws.Range("E30").Value = '  +3.16%  '

$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = '1.532'
$r.Style = "Normal"
$ws.Range("E31").Value = '  -0.79%  '

$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = '4.502'
$r.Style = "Normal"
$ws.Range("E32").Value = '  +0.49%  '

$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = '4.165'
$r.Style = "Normal"
$ws.Range("E33").Value = '  +2.29%  '

$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = '0.05415'
$r.Style = "Normal"
$ws.Range("E34").Value = '  -3.03%  '

$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = '1.238'
$r.Style = "Normal"
$ws.Range("E35").Value = '  -0.41%  '

$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = '0.7517'
$r.Style = "Normal"
$ws.Range("E36").Value = '  -0.42%  '

$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = '1.004'
$r.Style = "Normal"
$ws.Range("E37").Value = '  +0.09%  '

$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = '2.707'
$r.Style = "Normal"
$ws.Range("E38").Value = '  +0.94%  '

$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.01944'
$r.Style = "Normal"
$ws.Range("E39").Value = '  +0.56%  '

$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '2.769'
$r.Style = "Normal"
$ws.Range("E40").Value = '  -0.86%  '

$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.4460'
$r.Style = "Normal"
$ws.Range("E41").Value = '  -0.33%  '

$ws.Range("D42").Value = '1.098.53'
$ws.Range("E42").Value = '  -0.64%  '

$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = '6.074'
$r.Style = "Normal"
$ws.Range("E43").Value = '  +0.52%  '

$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '72.33'
$r.Style = "Normal"
$ws.Range("E44").Value = '  -3.08%  '

$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.8591'
$r.Style = "Normal"
$ws.Range("E45").Value = '  +0.97%  '

$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = '1.0000'
$r.Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '

$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = '7.722'
$r.Style = "Normal"
$ws.Range("E47").Value = '  +2.32%  '

$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '102.11'
$r.Style = "Normal"
$ws.Range("E48").Value = '  -0.44%  '

$ws.Range("E49").Value = '  -2.18%  '

$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = '3.030'
$r.Style = "Normal"
$ws.Range("E50").Value = '  +0.67%  '

$ws.Range("D51").Value = '2.028.18'
$ws.Range("E51").Value = '  -1.88%  '
